# fix: unique command names in XLSX - prefix protocol name to each step
#
# Every "protocol" worksheet (price1, price2, discount1, discount2, free1, ...,
# cumcontrol, dickpic, boosters) has a table with a header row ("Name", "Text",
# "Note", "*Guidelines") in row 1 and one row per command/step below it. The
# command names in column A used to collide across sheets (e.g. "Step1 Firmness"
# appeared in both discount1 and discount2), so we prefix each one with its
# sheet (protocol) name to make it unique, e.g. "discount2 Step1 Firmness".
#
# The first six "journey" worksheets (PeterJourney, MeetupRedirect, NRWaves,
# PersonalPeter, PositiveSpin, ReEngagement) are untouched.

$wb = $excel.ActiveWorkbook

$protocolSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $protocolSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.UsedRange.Rows.Count
    $firstRow = $ws.UsedRange.Row
    $lastDataRow = $firstRow + $lastRow - 1

    # Row 1 is the header ("Name"/"Text"/"Note"/"*Guidelines"); data rows start at 2.
    for ($r = 2; $r -le $lastDataRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $existing = $cell.Value2

        if ($existing -ne $null -and $existing -ne "") {
            $prefix = $sheetName + " "
            if ($existing.ToString().StartsWith($prefix) -eq $false) {
                $cell.Value2 = $prefix + $existing
            }
        }
    }
}
